$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.943.03'
$ws.Range("E2").Value = '  -0.51%  '

$ws.Range("D3").Value = '2.215.31'
$ws.Range("E3").Value = '  -1.26%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.86'
$ws.Range("E5").Value = '  -2.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  +0.61%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.12'
$ws.Range("E7").Value = '  -1.55%  '

$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("E9").Value = '  -2.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.82'
$ws.Range("E10").Value = '  +1.67%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0951'
$ws.Range("E11").Value = '  +0.93%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.07'
$ws.Range("E12").Value = '  -1.39%  '

$ws.Range("E13").Value = '  -0.40%  '

$ws.Range("D14").Value = '2.547.06'
$ws.Range("E14").Value = '  -1.20%  '

$ws.Range("E15").Value = '  -1.52%  '

$ws.Range("E16").Value = '  -1.72%  '

$ws.Range("D17").Value = '2.217.66'
$ws.Range("E17").Value = '  -0.86%  '

$ws.Range("D18").Value = '41.809.49'
$ws.Range("E18").Value = '  -0.54%  '

$ws.Range("E19").Value = '  +10.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '73.36'
$ws.Range("E20").Value = '  +1.71%  '

$ws.Range("E21").Value = '  +0.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.37'
$ws.Range("E22").Value = '  +16.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '229.25'
$ws.Range("E23").Value = '  -0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.08'
$ws.Range("E24").Value = '  -5.94%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.62'
$ws.Range("E25").Value = '  +1.08%  '

$ws.Range("E26").Value = '  +0.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.61'
$ws.Range("E27").Value = '  -1.01%  '

$ws.Range("E28").Value = '  -1.67%  '

$ws.Range("E29").Value = '  +1.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.53'
$ws.Range("E30").Value = '  -1.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.57'
$ws.Range("E31").Value = '  -0.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.63'
$ws.Range("E32").Value = '  +8.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0793'
$ws.Range("E33").Value = '  -3.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.124'
$ws.Range("E34").Value = '  -0.50%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '29.13'
$ws.Range("E35").Value = '  -6.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.111'
$ws.Range("E36").Value = '  -7.84%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.23'
$ws.Range("E37").Value = '  -5.30%  '

$ws.Range("E38").Value = '  -3.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.41'
$ws.Range("E39").Value = '  -2.14%  '

$ws.Range("E40").Value = '  +6.29%  '

$ws.Range("E41").Value = '  -3.25%  '

$ws.Range("E42").Value = '  -2.65%  '

$ws.Range("E43").Value = '  -4.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.69'
$ws.Range("E44").Value = '  +0.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.47'
$ws.Range("E45").Value = '  -3.09%  '

$ws.Range("E46").Value = '  -2.40%  '

$ws.Range("E47").Value = '  +4.15%  '

$ws.Range("E48").Value = '  +0.28%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.17'
$ws.Range("E49").Value = '  +0.07%  '

$ws.Range("E50").Value = '  -0.32%  '

$ws.Range("D51").Value = '2.421.25'
$ws.Range("E51").Value = '  -1.35%  '
